$p = $ppt.ActivePresentation

# Remove the last three slides:
#   16 "Параметры метода для API"
#   17 "Задание - часть 1"
#   18 "Задание – часть 2"
# Delete from the highest index down so earlier indices stay valid.
$p.Slides.Item(18).Delete()
$p.Slides.Item(17).Delete()
$p.Slides.Item(16).Delete()
